$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A61:B61").Copy()
$ws.Range("A62:B62").PasteSpecial()

$ws.Range("A62").Value = "17-11-2025"
$ws.Range("B62").Value = "The price of gold in India today is ₹12,497 per gram for 24 karat gold, ₹11,455 per gram for 22 karat gold and ₹9,373 per gram for 18 karat gold (also called 999 gold)."
